$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1178792975459519
$ws.Range("C2").Value = 0.9978466243962513
$ws.Range("D2").Value = 0.2488675586823569
$ws.Range("G2").Value = 0.2675185060000028
$ws.Range("H2").Value = 0.987

$ws.Range("B3").Value = 0.2448419781721476
$ws.Range("C3").Value = 0.9819101826999579
$ws.Range("D3").Value = 0.3809880346727105
$ws.Range("G3").Value = 0.2675185060000028
$ws.Range("H3").Value = 0.987
